$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (b0)
$ws.Range("C2").Value = 1009
$ws.Range("D2").Value = 25.28
$ws.Range("E2").Value = 2.007
$ws.Range("F2").Value = 958
$ws.Range("G2").Value = 1011
$ws.Range("H2").Value = 1055

# Row 3 (b1)
$ws.Range("C3").Value = -1.192
$ws.Range("D3").Value = 1.087
$ws.Range("E3").Value = 0.08632
$ws.Range("F3").Value = -3.161
$ws.Range("G3").Value = -1.276
$ws.Range("H3").Value = 1.006

# Row 4 (b2) - E4 also loses its scientific-notation number format
$ws.Range("C4").Value = -171.2
$ws.Range("D4").Value = 23.44
$ws.Range("E4").ClearFormats()
$ws.Range("E4").Value = 0.2167
$ws.Range("F4").Value = -217
$ws.Range("G4").Value = -171.3
$ws.Range("H4").Value = -125.1

# Row 5 (mu1)
$ws.Range("C5").Value = 980.2
$ws.Range("D5").Value = 3.457
$ws.Range("E5").Value = 0.07335
$ws.Range("F5").Value = 973.4
$ws.Range("G5").Value = 980.2
$ws.Range("H5").Value = 987.1

# Row 6 (mu2)
$ws.Range("C6").Value = 804.2
$ws.Range("D6").Value = 24.05
$ws.Range("E6").Value = 0.4661
$ws.Range("F6").Value = 757
$ws.Range("G6").Value = 804.2
$ws.Range("H6").Value = 851.3

# Row 7 (pred1)
$ws.Range("C7").Value = 975.7
$ws.Range("D7").Value = 429.9
$ws.Range("E7").Value = 4.547
$ws.Range("F7").Value = 142.8
$ws.Range("G7").Value = 973.6
$ws.Range("H7").Value = 1835

# Row 8 (pred2)
$ws.Range("C8").Value = 802.5
$ws.Range("D8").Value = 425.8
$ws.Range("E8").Value = 4.136
$ws.Range("F8").Value = -36.39
$ws.Range("G8").Value = 803.5
$ws.Range("H8").Value = 1629

# Row 9 (sigma)
$ws.Range("C9").Value = 426.5
$ws.Range("D9").Value = 2.453
$ws.Range("E9").Value = 0.05907
$ws.Range("F9").Value = 421.8
$ws.Range("G9").Value = 426.5
$ws.Range("H9").Value = 431.3

# Row 10 (tau)
$ws.Range("C10").Value = 0.000005498
$ws.Range("D10").Value = 0.00000006325
$ws.Range("E10").Value = 0.000000001524
$ws.Range("F10").Value = 0.000005376
$ws.Range("G10").Value = 0.000005498
$ws.Range("H10").Value = 0.000005622

$wb.Save()
